# Update countries & provincias Spain
# Refresh of the "Pais" COVID dashboard sheet: new scrape timestamp, a handful
# of countries whose case totals changed enough to swap rank with a
# neighbouring row, and the resulting data refresh for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) -------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 22 de Junio de 2020 a las 07:54"

# --- Row 44: Afganistan -----------------------------------------------------
$ws.Cells.Item(44, 2).Value = 29143
$ws.Cells.Item(44, 3).Value = 310
$ws.Cells.Item(44, 4).Value = 8841
$ws.Cells.Item(44, 5).Value = 19704
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 17
$ws.Cells.Item(44, 8).Value = 598

# --- Row 76: Uzbekistan ------------------------------------------------------
$ws.Cells.Item(76, 2).Value = 6358
$ws.Cells.Item(76, 3).Value = 43
$ws.Cells.Item(76, 4).Value = 4377
$ws.Cells.Item(76, 5).Value = 1962
$ws.Cells.Item(76, 6).Value = 0
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 19

# --- Rows 80/81: Haiti overtakes Republica de Macedonia --------------------
$ws.Cells.Item(80, 1).Value = "Haiti"
$ws.Cells.Item(80, 2).Value = 5211
$ws.Cells.Item(80, 3).Value = 134
$ws.Cells.Item(80, 4).Value = 24
$ws.Cells.Item(80, 5).Value = 5099
$ws.Cells.Item(80, 6).Value = 0
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 88

$ws.Cells.Item(81, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(81, 2).Value = 5106
$ws.Cells.Item(81, 3).Value = 0
$ws.Cells.Item(81, 4).Value = 1926
$ws.Cells.Item(81, 5).Value = 2942
$ws.Cells.Item(81, 6).Value = 0
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = 238

# --- Row 84: El Salvador -----------------------------------------------------
$ws.Cells.Item(84, 4).Value = 2542
$ws.Cells.Item(84, 5).Value = 1977
$ws.Cells.Item(84, 7).Value = 9
$ws.Cells.Item(84, 8).Value = 107

# --- Row 95: Tailandia --------------------------------------------------------
$ws.Cells.Item(95, 2).Value = 3151
$ws.Cells.Item(95, 3).Value = 3
$ws.Cells.Item(95, 4).Value = 3022
$ws.Cells.Item(95, 5).Value = 71

# --- Row 185: Butan ------------------------------------------------------------
$ws.Cells.Item(185, 4).Value = 32
$ws.Cells.Item(185, 5).Value = 36

# --- Rows 202/203: Dominica / Fiyi swap ranking (tied totals) ---------------
$ws.Cells.Item(202, 1).Value = "Dominica"
$ws.Cells.Item(203, 1).Value = "Fiyi"

# --- Rows 207/208: Islas Malvinas / Groenlandia swap ranking (tied totals) --
$ws.Cells.Item(207, 1).Value = "Islas Malvinas"
$ws.Cells.Item(208, 1).Value = "Groenlandia"

# --- Row 212: Sahara Occidental ------------------------------------------------
$ws.Cells.Item(212, 2).Value = 10
$ws.Cells.Item(212, 3).Value = 1
$ws.Cells.Item(212, 5).Value = 1

# --- Rows 213/214: Islas Virgenes Britanicas / Papua Nueva Guinea swap ------
$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 8).Value = 1

$ws.Cells.Item(214, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 8).Value = 0
